# Apply updated renewable cost figures (IRA eligibility assessment, wind and
# solar capex updates) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Wind CapEx rows (9-13): columns B,C,D,E,F ---
$ws.Range("B9").Value  = 1302
$ws.Range("C9").Value  = 1302
$ws.Range("D9").Value  = 1422
$ws.Range("E9").Value  = 1873
$ws.Range("F9").Value  = 1302

$ws.Range("B10").Value = 1211
$ws.Range("C10").Value = 1211
$ws.Range("D10").Value = 1301
$ws.Range("E10").Value = 1692
$ws.Range("F10").Value = 1211

$ws.Range("B11").Value = 1098
$ws.Range("C11").Value = 1098
$ws.Range("D11").Value = 1150
$ws.Range("E11").Value = 1467
$ws.Range("F11").Value = 1098

$ws.Range("B12").Value = 1044
$ws.Range("C12").Value = 1044
$ws.Range("D12").Value = 1093
$ws.Range("E12").Value = 1396
$ws.Range("F12").Value = 1044

$ws.Range("B13").Value = 882
$ws.Range("C13").Value = 882
$ws.Range("D13").Value = 924
$ws.Range("E13").Value = 1185
$ws.Range("F13").Value = 882

# --- Wind OpEx rows (15-19): columns B,C,D,E,F ---
$ws.Range("B15").Value = 28.9
$ws.Range("C15").Value = 28.9
$ws.Range("D15").Value = 27.2
$ws.Range("E15").Value = 28.9
$ws.Range("F15").Value = 28.9

$ws.Range("B16").Value = 27.5
$ws.Range("C16").Value = 27.5
$ws.Range("D16").Value = 25.6
$ws.Range("E16").Value = 27.5
$ws.Range("F16").Value = 27.5

$ws.Range("B17").Value = 25.8
$ws.Range("C17").Value = 25.8
$ws.Range("D17").Value = 23.5
$ws.Range("E17").Value = 25.8
$ws.Range("F17").Value = 25.8

$ws.Range("B18").Value = 24.9
$ws.Range("C18").Value = 24.9
$ws.Range("D18").Value = 22.7
$ws.Range("E18").Value = 24.9
$ws.Range("F18").Value = 24.9

$ws.Range("B19").Value = 22.3
$ws.Range("C19").Value = 22.3
$ws.Range("D19").Value = 20.3
$ws.Range("E19").Value = 22.3
$ws.Range("F19").Value = 22.3

# --- PV base installed cost rows (22-26): uniform across columns B:F ---
$ws.Range("B22:F22").Value = 1233
$ws.Range("B23:F23").Value = 1192
$ws.Range("B24:F24").Value = 991
$ws.Range("B25:F25").Value = 792
$ws.Range("B26:F26").Value = 604

# --- PV OpEx rows (28-32): uniform across columns B:F ---
$ws.Range("B28:F28").Value = 21.5
$ws.Range("B29:F29").Value = 19.6
$ws.Range("B30:F30").Value = 17.2
$ws.Range("B31:F31").Value = 14.8
$ws.Range("B32:F32").Value = 12.9

# Update the active selection to reflect where the author left off editing.
$ws.Range("D23").Select()
